# Ajustement du nombre d'équipe pour calcul heure fin du CLMI
# - 13 équipes de 6 coureurs
# - Annulation de la présence de l'Algérie

$wb = $excel.ActiveWorkbook
$wsEquipes = $wb.Worksheets.Item("EQUIPES")

# Remove the "Équipe nationale d'Algérie" row (row 6) from the EQUIPES sheet.
# Deleting the whole row shifts the rows below it up by one, and Excel
# automatically drops the now-unused shared string and reindexes the
# remaining shared strings / formulas accordingly.
$wsEquipes.Rows.Item(6).Delete()

# The active sheet moves from ADMIN back to EQUIPES, with a fresh selection.
$wsEquipes.Activate()
$wsEquipes.Range("A24").Select()
